# Adjust api response formatting
# Fix various typos in comments and documentation.
# This resolves issue #17.
#
# Append a new row (row 44) of decoded packet data to each of the four
# worksheets, duplicating the structure/format of the preceding row (43)
# with updated values.

$wb = $excel.ActiveWorkbook

function Add-DataRow($ws, $row, $timeValue, $colB, $colC, $colD, $colE, $colF, $colG, $colH, $colI) {
    $prevRow = $row - 1

    $ws.Cells.Item($row, 1).Value = $timeValue
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $colD
    $ws.Cells.Item($row, 5).Value = $colE

    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
}

$timeValue = 45830.49516203703

# Sheet 1: FE_LFT_#1
$ws1 = $wb.Worksheets.Item(1)
$g1 = [double]"7.598631275147109e+23"
Add-DataRow $ws1 44 $timeValue "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x60" "0xf" 380 $g1 352 15

# Sheet 2: FE_LFT_#2
$ws2 = $wb.Worksheets.Item(2)
$g2 = [double]"5.68432987514711e+23"
Add-DataRow $ws2 44 $timeValue "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x74" "0xe" 400 $g2 372 14

# Sheet 3: FE_PLT_#1
$ws3 = $wb.Worksheets.Item(3)
$g3 = [double]"5.68631262647114e+23"
Add-DataRow $ws3 44 $timeValue "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x6A" "0x3" 110 $g3 106 3

# Sheet 4: FE_PLT_#2
$ws4 = $wb.Worksheets.Item(4)
$g4 = [double]"9.85046333984776e+23"
Add-DataRow $ws4 44 $timeValue "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x6A" "0x3" 110 $g4 106 3
